$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Thomas Hex" -> "Matthies Hex" (text stays in same data row, no position shift) ---
$ws.Range("B9").Value() = "Matthies Hex"

# --- Insert two new rows for "Holden" and "Rizzie Spiral" ahead of the old row 4 block ---
$ws.Rows("4:5").Insert()

# Copy formatting (bold / border / alignment) from row 3's index cell onto the two new index cells
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4: new "Holden" measurement
$ws.Range("A4").Value() = 2
$ws.Range("B4").Value() = "Holden"
$ws.Range("C4").Value() = 4.723790133246625
$ws.Range("D4").Value() = 0.8243751953151571
$ws.Range("E4").Value() = 0.9950484231798352
$ws.Range("F4").Value() = 0.9719426349527742
$ws.Range("G4").Value() = 0.9460216726130078
$ws.Range("H4").Value() = 0.6979339306661129
$ws.Range("I4").Value() = 0.6979339306661129
$ws.Range("J4").Value() = 4.723790133246625
$ws.Range("K4").Value() = 4.723790133246625
$ws.Range("L4").Value() = 0.9460216726130078
$ws.Range("M4").Value() = 0.8219778016395604
$ws.Range("N4").Value() = 0.8219778016395604
$ws.Range("O4").Value() = 0.8227769328647593
$ws.Range("P4").Value() = 2.122581912175248
$ws.Range("Q4").Value() = 2.122581912175248
$ws.Range("R4").Value() = 2.772883967443093
$ws.Range("S4").Value() = 2.772883967443093
$ws.Range("T4").Value() = 1.526518664995585

# Row 5: new "Rizzie Spiral" measurement
$ws.Range("A5").Value() = 3
$ws.Range("B5").Value() = "Rizzie Spiral"
$ws.Range("C5").Value() = 0.6933168727090753
$ws.Range("D5").Value() = 0.002216029919211659
$ws.Range("E5").Value() = 11.85567036476542
$ws.Range("F5").Value() = -0.00725402253604594
$ws.Range("G5").Value() = 3.957814197546023
$ws.Range("H5").Value() = 8.979342344833363
$ws.Range("I5").Value() = 8.979342344833363
$ws.Range("J5").Value() = 0.6933168727090753
$ws.Range("K5").Value() = 0.6933168727090753
$ws.Range("L5").Value() = 3.957814197546023
$ws.Range("M5").Value() = 6.468578271189693
$ws.Range("N5").Value() = 6.468578271189693
$ws.Range("O5").Value() = 4.313124190766199
$ws.Range("P5").Value() = 4.54349113836282
$ws.Range("Q5").Value() = 4.54349113836282
$ws.Range("R5").Value() = 3.580947571949384
$ws.Range("S5").Value() = 3.580947571949384
$ws.Range("T5").Value() = 4.246850964539507
